$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (预测值) updates for rows 2-24
$ws.Range("C2").Value = -1.6098
$ws.Range("C3").Value = -2.2803
$ws.Range("C4").Value = -2.6555
$ws.Range("C5").Value = -2.7659
$ws.Range("C6").Value = -2.7906
$ws.Range("C7").Value = -2.5474
$ws.Range("C8").Value = -2.7462
$ws.Range("C9").Value = -2.3399
$ws.Range("C10").Value = -1.3707
$ws.Range("C11").Value = -0.7468
$ws.Range("C12").Value = -0.243
$ws.Range("C13").Value = 0.2667
$ws.Range("C14").Value = -0.2475
$ws.Range("C15").Value = -0.7000999999999999
$ws.Range("C16").Value = -1.0698
$ws.Range("C17").Value = -0.893
$ws.Range("C18").Value = -0.5856
$ws.Range("C19").Value = -1.3812
$ws.Range("C20").Value = -0.996
$ws.Range("C21").Value = -0.1794
$ws.Range("C22").Value = -0.1369
$ws.Range("C23").Value = -0.7907999999999999
$ws.Range("C24").Value = -0.9343

# Column B (真实值) update for row 3
$ws.Range("B3").Value = -5.2367
